$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update column F ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 0
$wsExhibit.Range("F3").Value = 0
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F9").Value = 0
$wsExhibit.Range("F10").Value = 515

# Sheet "全部类型" (All types) - update column F ("想去人数") values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 120
$wsAll.Range("F4").Value = 0
$wsAll.Range("F5").Value = 19
$wsAll.Range("F6").Value = 0
$wsAll.Range("F7").Value = 0
$wsAll.Range("F10").Value = 515
